$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,7).Value = 4.882566666666666
$ws.Cells.Item(2,8).Value = 14.6477
$ws.Cells.Item(2,9).Value = 0.1179541876619502
$ws.Cells.Item(2,10).Value = 0.1179541876619503
$ws.Cells.Item(2,13).Value = 0.667106
$ws.Cells.Item(2,14).Value = 2.001318
$ws.Cells.Item(2,15).Value = 0.003817114239487378
$ws.Cells.Item(2,16).Value = 0.003817114239487378
$ws.Cells.Item(2,17).Value = 3.257189518733333
$ws.Cells.Item(2,18).Value = 29.3147056686
$ws.Cells.Item(2,19).Value = 0.0004502446093315967
$ws.Cells.Item(2,20).Value = 0.0004502446093315967

# Row 3
$ws.Cells.Item(3,7).Value = 4.882566666666666
$ws.Cells.Item(3,8).Value = 14.6477
$ws.Cells.Item(3,9).Value = 0.1179541876619502
$ws.Cells.Item(3,10).Value = 0.1179541876619503
$ws.Cells.Item(3,15).Value = 0.9945745510447523
$ws.Cells.Item(3,16).Value = 0.9945745510447522
$ws.Cells.Item(3,17).Value = 848.6824339050777
$ws.Cells.Item(3,18).Value = 7638.1419051457
$ws.Cells.Item(3,19).Value = 0.1173142332377326
$ws.Cells.Item(3,20).Value = 0.1173142332377326

# Row 4
$ws.Cells.Item(4,7).Value = 4.882566666666666
$ws.Cells.Item(4,8).Value = 14.6477
$ws.Cells.Item(4,9).Value = 0.1179541876619502
$ws.Cells.Item(4,10).Value = 0.1179541876619503
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.260372
$ws.Cells.Item(4,14).Value = 0.7811159999999999
$ws.Cells.Item(4,15).Value = 0.001489822709979835
$ws.Cells.Item(4,16).Value = 0.001489822709979834
$ws.Cells.Item(4,17).Value = 1.271283648133333
$ws.Cells.Item(4,18).Value = 11.4415528332
$ws.Cells.Item(4,19).Value = 0.0001757308275159967
$ws.Cells.Item(4,20).Value = 0.0001757308275159967

# Row 5
$ws.Cells.Item(5,7).Value = 4.882566666666666
$ws.Cells.Item(5,8).Value = 14.6477
$ws.Cells.Item(5,9).Value = 0.1179541876619502
$ws.Cells.Item(5,10).Value = 0.1179541876619503
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.020712
$ws.Cells.Item(5,14).Value = 0.062136
$ws.Cells.Item(5,15).Value = 0.0001185120057805845
$ws.Cells.Item(5,16).Value = 0.0001185120057805844
$ws.Cells.Item(5,17).Value = 0.1011277208
$ws.Cells.Item(5,18).Value = 0.9101494872
$ws.Cells.Item(5,19).Value = 0.00001397898737003719
$ws.Cells.Item(5,20).Value = 0.00001397898737003719

# Row 6
$ws.Cells.Item(6,9).Value = 0.0529488011407969
$ws.Cells.Item(6,10).Value = 0.0529488011407969
$ws.Cells.Item(6,13).Value = 0.667106
$ws.Cells.Item(6,14).Value = 2.001318
$ws.Cells.Item(6,15).Value = 0.003817114239487378
$ws.Cells.Item(6,16).Value = 0.003817114239487378
$ws.Cells.Item(6,17).Value = 1.462129353131333
$ws.Cells.Item(6,18).Value = 13.159164178182
$ws.Cells.Item(6,19).Value = 0.0002021116227983214
$ws.Cells.Item(6,20).Value = 0.0002021116227983213

# Row 7
$ws.Cells.Item(7,9).Value = 0.0529488011407969
$ws.Cells.Item(7,10).Value = 0.0529488011407969
$ws.Cells.Item(7,15).Value = 0.9945745510447523
$ws.Cells.Item(7,16).Value = 0.9945745510447522
$ws.Cells.Item(7,19).Value = 0.05266153012296595
$ws.Cells.Item(7,20).Value = 0.05266153012296594

# Row 8
$ws.Cells.Item(8,9).Value = 0.0529488011407969
$ws.Cells.Item(8,10).Value = 0.0529488011407969
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.260372
$ws.Cells.Item(8,14).Value = 0.7811159999999999
$ws.Cells.Item(8,15).Value = 0.001489822709979835
$ws.Cells.Item(8,16).Value = 0.001489822709979834
$ws.Cells.Item(8,17).Value = 0.5706702442093333
$ws.Cells.Item(8,18).Value = 5.136032197883999
$ws.Cells.Item(8,19).Value = 0.0000788843264057654
$ws.Cells.Item(8,20).Value = 0.00007888432640576538

# Row 9
$ws.Cells.Item(9,9).Value = 0.0529488011407969
$ws.Cells.Item(9,10).Value = 0.0529488011407969
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.020712
$ws.Cells.Item(9,14).Value = 0.062136
$ws.Cells.Item(9,15).Value = 0.0001185120057805845
$ws.Cells.Item(9,16).Value = 0.0001185120057805844
$ws.Cells.Item(9,17).Value = 0.04539551909599999
$ws.Cells.Item(9,18).Value = 0.4085596718639999
$ws.Cells.Item(9,19).Value = 0.000006275068626873139
$ws.Cells.Item(9,20).Value = 0.000006275068626873138

# Row 10
$ws.Cells.Item(10,7).Value = 25.81983966666667
$ws.Cells.Item(10,8).Value = 77.459519
$ws.Cells.Item(10,9).Value = 0.6237617264369424
$ws.Cells.Item(10,10).Value = 0.6237617264369424
$ws.Cells.Item(10,13).Value = 0.667106
$ws.Cells.Item(10,14).Value = 2.001318
$ws.Cells.Item(10,15).Value = 0.003817114239487378
$ws.Cells.Item(10,16).Value = 0.003817114239487378
$ws.Cells.Item(10,17).Value = 17.22456996067133
$ws.Cells.Item(10,18).Value = 155.021129646042
$ws.Cells.Item(10,19).Value = 0.002380969768029683
$ws.Cells.Item(10,20).Value = 0.002380969768029683

# Row 11
$ws.Cells.Item(11,7).Value = 25.81983966666667
$ws.Cells.Item(11,8).Value = 77.459519
$ws.Cells.Item(11,9).Value = 0.6237617264369424
$ws.Cells.Item(11,10).Value = 0.6237617264369424
$ws.Cells.Item(11,15).Value = 0.9945745510447523
$ws.Cells.Item(11,16).Value = 0.9945745510447522
$ws.Cells.Item(11,17).Value = 4487.976481907508
$ws.Cells.Item(11,18).Value = 40391.78833716758
$ws.Cells.Item(11,19).Value = 0.6203775390299217
$ws.Cells.Item(11,20).Value = 0.6203775390299215

# Row 12
$ws.Cells.Item(12,7).Value = 25.81983966666667
$ws.Cells.Item(12,8).Value = 77.459519
$ws.Cells.Item(12,9).Value = 0.6237617264369424
$ws.Cells.Item(12,10).Value = 0.6237617264369424
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.260372
$ws.Cells.Item(12,14).Value = 0.7811159999999999
$ws.Cells.Item(12,15).Value = 0.001489822709979835
$ws.Cells.Item(12,16).Value = 0.001489822709979834
$ws.Cells.Item(12,17).Value = 6.722763293689333
$ws.Cells.Item(12,18).Value = 60.504869643204
$ws.Cells.Item(12,19).Value = 0.0009292943856619859
$ws.Cells.Item(12,20).Value = 0.0009292943856619856

# Row 13
$ws.Cells.Item(13,7).Value = 25.81983966666667
$ws.Cells.Item(13,8).Value = 77.459519
$ws.Cells.Item(13,9).Value = 0.6237617264369424
$ws.Cells.Item(13,10).Value = 0.6237617264369424
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.020712
$ws.Cells.Item(13,14).Value = 0.062136
$ws.Cells.Item(13,15).Value = 0.0001185120057805845
$ws.Cells.Item(13,16).Value = 0.0001185120057805844
$ws.Cells.Item(13,17).Value = 0.5347805191759999
$ws.Cells.Item(13,18).Value = 4.813024672584
$ws.Cells.Item(13,19).Value = 0.00007392325332920226
$ws.Cells.Item(13,20).Value = 0.00007392325332920225

# Row 14
$ws.Cells.Item(14,7).Value = 8.499598333333333
$ws.Cells.Item(14,8).Value = 25.498795
$ws.Cells.Item(14,9).Value = 0.2053352847603104
$ws.Cells.Item(14,10).Value = 0.2053352847603104
$ws.Cells.Item(14,13).Value = 0.667106
$ws.Cells.Item(14,14).Value = 2.001318
$ws.Cells.Item(14,15).Value = 0.003817114239487378
$ws.Cells.Item(14,16).Value = 0.003817114239487378
$ws.Cells.Item(14,17).Value = 5.670133045756667
$ws.Cells.Item(14,18).Value = 51.03119741181
$ws.Cells.Item(14,19).Value = 0.0007837882393277765
$ws.Cells.Item(14,20).Value = 0.0007837882393277765

# Row 15
$ws.Cells.Item(15,7).Value = 8.499598333333333
$ws.Cells.Item(15,8).Value = 25.498795
$ws.Cells.Item(15,9).Value = 0.2053352847603104
$ws.Cells.Item(15,10).Value = 0.2053352847603104
$ws.Cells.Item(15,15).Value = 0.9945745510447523
$ws.Cells.Item(15,16).Value = 0.9945745510447522
$ws.Cells.Item(15,17).Value = 1477.39094890301
$ws.Cells.Item(15,18).Value = 13296.5185401271
$ws.Cells.Item(15,19).Value = 0.2042212486541321
$ws.Cells.Item(15,20).Value = 0.2042212486541321

# Row 16
$ws.Cells.Item(16,7).Value = 8.499598333333333
$ws.Cells.Item(16,8).Value = 25.498795
$ws.Cells.Item(16,9).Value = 0.2053352847603104
$ws.Cells.Item(16,10).Value = 0.2053352847603104
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.260372
$ws.Cells.Item(16,14).Value = 0.7811159999999999
$ws.Cells.Item(16,15).Value = 0.001489822709979835
$ws.Cells.Item(16,16).Value = 0.001489822709979834
$ws.Cells.Item(16,17).Value = 2.213057417246667
$ws.Cells.Item(16,18).Value = 19.91751675522
$ws.Cells.Item(16,19).Value = 0.0003059131703960867
$ws.Cells.Item(16,20).Value = 0.0003059131703960867

# Row 17
$ws.Cells.Item(17,7).Value = 8.499598333333333
$ws.Cells.Item(17,8).Value = 25.498795
$ws.Cells.Item(17,9).Value = 0.2053352847603104
$ws.Cells.Item(17,10).Value = 0.2053352847603104
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.020712
$ws.Cells.Item(17,14).Value = 0.062136
$ws.Cells.Item(17,15).Value = 0.0001185120057805845
$ws.Cells.Item(17,16).Value = 0.0001185120057805844
$ws.Cells.Item(17,17).Value = 0.17604368068
$ws.Cells.Item(17,18).Value = 1.58439312612
$ws.Cells.Item(17,19).Value = 0.00002433469645447186
$ws.Cells.Item(17,20).Value = 0.00002433469645447186

